$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 3.1851815
$ws.Range("N2").Value = 6.370363
$ws.Range("O2").Value = 0.4406530230187619
$ws.Range("P2").Value = 0.3851702893788179
$ws.Range("Q2").Value = 1.189168401936
$ws.Range("R2").Value = 7.135010411615999
$ws.Range("S2").Value = 0.3230233816467307
$ws.Range("T2").Value = 0.3099293428980676
$ws.Range("O3").Value = 0.2827048402157753
$ws.Range("P3").Value = 0.3706641033643825
$ws.Range("S3").Value = 0.2072385045013301
$ws.Range("T3").Value = 0.298256862378706
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.003928
$ws.Range("N4").Value = 0.011784
$ws.Range("O4").Value = 0.0005434180358066555
$ws.Range("P4").Value = 0.0007124942001013113
$ws.Range("Q4").Value = 0.001466495232
$ws.Range("R4").Value = 0.013198457088
$ws.Range("S4").Value = 0.000398355900003927
$ws.Range("T4").Value = 0.0005733122863973101
$ws.Range("M5").Value = 1.9606995
$ws.Range("N5").Value = 3.921399
$ws.Range("O5").Value = 0.2712524111754306
$ws.Range("P5").Value = 0.2370989514411984
$ws.Range("Q5").Value = 0.7320153941279999
$ws.Range("R5").Value = 4.392092364768
$ws.Range("S5").Value = 0.1988432316598141
$ws.Range("T5").Value = 0.1907829452279469
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02851766666666667
$ws.Range("N6").Value = 0.085553
$ws.Range("O6").Value = 0.003945268433245655
$ws.Range("P6").Value = 0.005172778029639129
$ws.Range("Q6").Value = 0.010646899744
$ws.Range("R6").Value = 0.095822097696
$ws.Range("S6").Value = 0.002892103047609976
$ws.Range("T6").Value = 0.004162303635280811
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.006513000000000001
$ws.Range("N7").Value = 0.019539
$ws.Range("O7").Value = 0.0009010391209798237
$ws.Range("P7").Value = 0.001181383585860448
$ws.Range("Q7").Value = 0.002431589472
$ws.Range("R7").Value = 0.021884305248
$ws.Range("S7").Value = 0.0006605122140340063
$ws.Range("T7").Value = 0.0009506066500269046
$ws.Range("M8").Value = 3.1851815
$ws.Range("N8").Value = 6.370363
$ws.Range("O8").Value = 0.4406530230187619
$ws.Range("P8").Value = 0.3851702893788179
$ws.Range("Q8").Value = 0.433038165651
$ws.Range("R8").Value = 1.732152662604
$ws.Range("S8").Value = 0.1176296413720312
$ws.Range("T8").Value = 0.07524094648075033
$ws.Range("O9").Value = 0.2827048402157753
$ws.Range("P9").Value = 0.3706641033643825
$ws.Range("S9").Value = 0.0754663357144452
$ws.Range("T9").Value = 0.07240724098567648
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.003928
$ws.Range("N10").Value = 0.011784
$ws.Range("O10").Value = 0.0005434180358066555
$ws.Range("P10").Value = 0.0007124942001013113
$ws.Range("Q10").Value = 0.000534027312
$ws.Range("R10").Value = 0.003204163871999999
$ws.Range("S10").Value = 0.0001450621358027286
$ws.Range("T10").Value = 0.0001391819137040012
$ws.Range("M11").Value = 1.9606995
$ws.Range("N11").Value = 3.921399
$ws.Range("O11").Value = 0.2712524111754306
$ws.Range("P11").Value = 0.2370989514411984
$ws.Range("Q11").Value = 0.266564939823
$ws.Range("R11").Value = 1.066259759292
$ws.Range("S11").Value = 0.0724091795156166
$ws.Range("T11").Value = 0.04631600621325157
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02851766666666667
$ws.Range("N12").Value = 0.085553
$ws.Range("O12").Value = 0.003945268433245655
$ws.Range("P12").Value = 0.005172778029639129
$ws.Range("Q12").Value = 0.003877090854
$ws.Range("R12").Value = 0.023262545124
$ws.Range("S12").Value = 0.001053165385635679
$ws.Range("T12").Value = 0.001010474394358317
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.006513000000000001
$ws.Range("N13").Value = 0.019539
$ws.Range("O13").Value = 0.0009010391209798237
$ws.Range("P13").Value = 0.001181383585860448
$ws.Range("Q13").Value = 0.0008854684020000001
$ws.Range("R13").Value = 0.005312810411999999
$ws.Range("S13").Value = 0.0002405269069458175
$ws.Range("T13").Value = 0.0002307769358335437
